$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1: "idReponse" / "idQuestionChoixReponse") is removed.
# Deleting the entire row shifts every subsequent data row up by one and
# drops the now out-of-range last row, matching the target sheet exactly.
$ws.Rows(1).Select() | Out-Null
$ws.Rows(1).Delete()
